$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 2 ("R Rohit"). This shifts
# "R Rohit" down to row 3 and "Jagadeesh Latti" down to row 4, preserving
# their existing cell contents.
$ws.Rows.Item(2).Insert()

# --- Row 2: new attendee "Varun" ---
$ws.Range("A2").Value = "Varun"

# Excel's native Insert() leaves the new row's cells completely blank;
# the workbook's convention (used throughout this sheet) is that an
# absent attendance mark is stored as an explicit empty-text cell rather
# than a truly blank cell. A leading single-quote forces text type with
# empty content, then ClearFormats() drops the quote-prefix style it
# introduces so the cell's style stays the default (style 0).
$cols = @("B","C","D","E","F","G","H","I","J","K")
foreach ($col in $cols) {
    $cell = $ws.Range($col + "2")
    $cell.Value = "'"
    $cell.ClearFormats()
}

# Varun attended (P) on Feb-27 and Feb-28
$ws.Range("L2").Value = "P"
$ws.Range("M2").Value = "P"

# --- Row 3: "R Rohit" (shifted down from row 2) ---
# Only Feb-28 attendance is recorded for him; the rest stay as the
# empty-text cells that were already shifted down with the row.
$ws.Range("M3").Value = "P"

# --- Row 4: "Jagadeesh Latti" (shifted down from row 3) ---
# Fill in the newly recorded attendance codes.
$ws.Range("C4").Value = "U"
$ws.Range("D4").Value = "W"
$ws.Range("G4").Value = "U"
$ws.Range("H4").Value = "W"
$ws.Range("J4").Value = "W"

# Feb-27/Feb-28 are no longer marked for him - reset to the sheet's
# empty-text convention (same trick as above).
foreach ($col in @("L","M")) {
    $cell = $ws.Range($col + "4")
    $cell.Value = "'"
    $cell.ClearFormats()
}
